$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 20081958
$ws.Range("I132").Value = 21363220
$ws.Range("K132").Value = 64089660
$ws.Range("M132").Value = -64087130
$ws.Range("H137").Value = 3367.587
$ws.Range("I137").Value = 2892.8684
$ws.Range("J137").Value = 5622.5
$ws.Range("K137").Value = 8678.6052
$ws.Range("L137").Value = 16867.5
$ws.Range("M137").Value = -6128.6052
$ws.Range("N137").Value = -21967.5
$ws.Range("H138").Value = 2577.5276
$ws.Range("I138").Value = 1476.4375
$ws.Range("J138").Value = 3174.7288
$ws.Range("K138").Value = 4429.3125
$ws.Range("L138").Value = 9524.186399999999
$ws.Range("M138").Value = 710.6875
$ws.Range("N138").Value = -19804.1864

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9126.696
$ws.Range("I32").Value = 5952.357
$ws.Range("J32").Value = 14513.454
$ws.Range("K32").Value = 5952.357
$ws.Range("L32").Value = 14513.454
$ws.Range("M32").Value = -5665.357
$ws.Range("N32").Value = -15087.454
$ws.Range("H74").Value = 1241.4584
$ws.Range("I74").Value = 904.27905
$ws.Range("K74").Value = 904.27905
$ws.Range("M74").Value = -30.27904999999998
$ws.Range("H77").Value = 1241.4584
$ws.Range("I77").Value = 904.27905
$ws.Range("K77").Value = 4521.39525
$ws.Range("M77").Value = -153.3952499999996
$ws.Range("H97").Value = 759.13043
$ws.Range("I97").Value = 613.6842
$ws.Range("J97").Value = 1450
$ws.Range("K97").Value = 613.6842
$ws.Range("L97").Value = 1450
$ws.Range("M97").Value = -117.6842
$ws.Range("N97").Value = -2442
$ws.Range("H132").Value = 2716.3389
$ws.Range("I132").Value = 2063.149
$ws.Range("K132").Value = 6189.447
$ws.Range("M132").Value = -3659.447

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2193.3157
$ws.Range("I86").Value = 2162.5715
$ws.Range("J86").Value = 2279.4
$ws.Range("K86").Value = 2162.5715
$ws.Range("L86").Value = 2279.4
$ws.Range("M86").Value = -1039.5715
$ws.Range("N86").Value = -4525.4
$ws.Range("H89").Value = 2193.3157
$ws.Range("I89").Value = 2162.5715
$ws.Range("J89").Value = 2279.4
$ws.Range("K89").Value = 10812.8575
$ws.Range("L89").Value = 11397
$ws.Range("M89").Value = -5196.8575
$ws.Range("N89").Value = -22629
$ws.Range("H105").Value = 2871.818
$ws.Range("I105").Value = 2871.818
$ws.Range("K105").Value = 2871.818
$ws.Range("M105").Value = -1124.818

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3079.3784
$ws.Range("I31").Value = 1067.2273
$ws.Range("J31").Value = 6030.533
$ws.Range("K31").Value = 1067.2273
$ws.Range("L31").Value = 6030.533
$ws.Range("M31").Value = -772.2273
$ws.Range("N31").Value = -6620.533
$ws.Range("H34").Value = 3079.3784
$ws.Range("I34").Value = 1067.2273
$ws.Range("J34").Value = 6030.533
$ws.Range("K34").Value = 1067.2273
$ws.Range("L34").Value = 6030.533
$ws.Range("M34").Value = -865.2273
$ws.Range("N34").Value = -6434.533
$ws.Range("H58").Value = 1990.6377
$ws.Range("I58").Value = 1644.5938
$ws.Range("J58").Value = 6420
$ws.Range("K58").Value = 1644.5938
$ws.Range("L58").Value = 6420
$ws.Range("M58").Value = -1441.5938
$ws.Range("N58").Value = -6826
$ws.Range("H136").Value = 1990.6377
$ws.Range("I136").Value = 1644.5938
$ws.Range("J136").Value = 6420
$ws.Range("K136").Value = 4933.7814
$ws.Range("L136").Value = 19260
$ws.Range("M136").Value = -2383.7814
$ws.Range("N136").Value = -24360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 398.1111
$ws.Range("I97").Value = 270.75
$ws.Range("K97").Value = 812.25
$ws.Range("M97").Value = -316.25
$ws.Range("H138").Value = 1753.8889
$ws.Range("I138").Value = 1204.6666
$ws.Range("J138").Value = 4500
$ws.Range("K138").Value = 3613.9998
$ws.Range("L138").Value = 13500
$ws.Range("M138").Value = 1526.0002
$ws.Range("N138").Value = -23780

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 50002244
$ws.Range("I80").Value = 250000000
$ws.Range("J80").Value = 2805.5
$ws.Range("K80").Value = 250000000
$ws.Range("L80").Value = 2805.5
$ws.Range("M80").Value = -249999002
$ws.Range("N80").Value = -4801.5
$ws.Range("H83").Value = 50002244
$ws.Range("I83").Value = 250000000
$ws.Range("J83").Value = 2805.5
$ws.Range("K83").Value = 1250000000
$ws.Range("L83").Value = 14027.5
$ws.Range("M83").Value = -1249995008
$ws.Range("N83").Value = -24011.5
$ws.Range("H122").Value = 2736.8
$ws.Range("I122").Value = 2418.12
$ws.Range("J122").Value = 3533.5
$ws.Range("K122").Value = 7254.36
$ws.Range("L122").Value = 10600.5
$ws.Range("M122").Value = -4804.36
$ws.Range("N122").Value = -15500.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2016.9231
$ws.Range("J16").Value = 2130
$ws.Range("L16").Value = 2130
$ws.Range("N16").Value = -2470
$ws.Range("H46").Value = 2613.1333
$ws.Range("I46").Value = 3799.3333
$ws.Range("J46").Value = 2316.5833
$ws.Range("K46").Value = 3799.3333
$ws.Range("L46").Value = 2316.5833
$ws.Range("M46").Value = -3611.3333
$ws.Range("N46").Value = -2692.5833
$ws.Range("H82").Value = 996.913
$ws.Range("I82").Value = 785.7368
$ws.Range("J82").Value = 2000
$ws.Range("K82").Value = 785.7368
$ws.Range("L82").Value = 2000
$ws.Range("M82").Value = -424.7368
$ws.Range("N82").Value = -2722
$ws.Range("H85").Value = 996.913
$ws.Range("I85").Value = 785.7368
$ws.Range("J85").Value = 2000
$ws.Range("K85").Value = 785.7368
$ws.Range("L85").Value = 2000
$ws.Range("M85").Value = 462.2632
$ws.Range("N85").Value = -4496
$ws.Range("H132").Value = 4303.7334
$ws.Range("I132").Value = 1253.4193
$ws.Range("J132").Value = 7564.4136
$ws.Range("K132").Value = 3760.2579
$ws.Range("L132").Value = 22693.2408
$ws.Range("M132").Value = -1230.2579
$ws.Range("N132").Value = -27753.2408

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 22960504
$ws.Range("I81").Value = 32144106
$ws.Range("J81").Value = 1495
$ws.Range("K81").Value = 64288212
$ws.Range("L81").Value = 2990
$ws.Range("M81").Value = -64287151
$ws.Range("N81").Value = -5112
$ws.Range("H84").Value = 22960504
$ws.Range("I84").Value = 32144106
$ws.Range("J84").Value = 1495
$ws.Range("K84").Value = 321441060
$ws.Range("L84").Value = 14950
$ws.Range("M84").Value = -321435756
$ws.Range("N84").Value = -25558
$ws.Range("H113").Value = 347.0476
$ws.Range("I113").Value = 327.7857
$ws.Range("J113").Value = 385.57144
$ws.Range("K113").Value = 983.3571000000001
$ws.Range("L113").Value = 1156.71432
$ws.Range("M113").Value = 1186.6429
$ws.Range("N113").Value = -5496.71432
$ws.Range("H126").Value = 6043.4595
$ws.Range("I126").Value = 2792.52
$ws.Range("J126").Value = 12816.25
$ws.Range("K126").Value = 8377.559999999999
$ws.Range("L126").Value = 38448.75
$ws.Range("M126").Value = -5907.559999999999
$ws.Range("N126").Value = -43388.75
$ws.Range("H132").Value = 7753867.5
$ws.Range("I132").Value = 650.7059
$ws.Range("J132").Value = 12823278
$ws.Range("K132").Value = 1952.1177
$ws.Range("L132").Value = 38469834
$ws.Range("M132").Value = 577.8822999999998
$ws.Range("N132").Value = -38474894
$ws.Range("H136").Value = 3372.9143
$ws.Range("I136").Value = 1019.6842
$ws.Range("K136").Value = 3059.0526
$ws.Range("M136").Value = -509.0526

Write-Output "applied 193 cell updates across 8 sheets"